function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$wb = $excel.ActiveWorkbook

# --- Insert a new worksheet "ARMS" immediately before "T2A" ---
$t2a = $wb.Worksheets.Item("T2A")
$ws = $wb.Worksheets.Add($t2a)
$ws.Name = "ARMS"

$blue  = RGBVal 0x38 0x77 0xA6   # header border - top/right (and left for col A)
$grey  = RGBVal 0xA5 0xA5 0xB1   # header border - bottom
$white = 16777215                # header font colour / fill pattern colour
$fillBlue = RGBVal 0x0B 0x64 0xA0  # header fill colour

# xlEdgeLeft=7, xlEdgeTop=8, xlEdgeBottom=9, xlEdgeRight=10, xlContinuous=1
function Set-HeaderCell($cell, $text, $wrap, $numFmt, $edges) {
    $cell.Value = $text

    $f = $cell.Font
    $f.Name = "Arial"
    $f.Bold = $true
    $f.Size = 9
    $f.Color = $white

    $cell.Interior.Color = $fillBlue
    $cell.Interior.PatternColor = $white

    $cell.HorizontalAlignment = -4131
    if ($wrap) {
        $cell.WrapText = $true
    } else {
        $cell.WrapText = $false
    }

    if ($numFmt) {
        $cell.NumberFormat = $numFmt
    }

    foreach ($e in $edges.Keys) {
        $cell.Borders.Item($e).LineStyle = 1
        $cell.Borders.Item($e).Color = $edges[$e]
    }
}

Set-HeaderCell $ws.Range("A1") "Assessment_Staff_Name"     $true  $null @{7=$blue; 8=$blue; 9=$grey; 10=$blue}
Set-HeaderCell $ws.Range("B1") "Assessment_Staff_Key"      $true  $null @{8=$blue; 9=$grey; 10=$blue}
Set-HeaderCell $ws.Range("C1") "Assessment_Staff_Grade"    $false "@"  @{8=$blue; 9=$grey; 10=$blue}
Set-HeaderCell $ws.Range("D1") "Assessmentent_Team_Key"    $true  $null @{8=$blue; 9=$grey; 10=$blue}
Set-HeaderCell $ws.Range("E1") "Assessment_Provider_Code"  $true  $null @{8=$blue; 9=$grey; 10=$blue}
Set-HeaderCell $ws.Range("F1") "CRN"                       $true  $null @{10=$blue}
Set-HeaderCell $ws.Range("G1") "Disposal_or_Release_Date"  $true  $null @{10=$blue}
Set-HeaderCell $ws.Range("H1") "Sentence_Type"             $true  $null @{10=$blue}
Set-HeaderCell $ws.Range("I1") "SO_Registration_Date"      $true  $null @{10=$blue}

# --- Row height for header row ---
$ws.Rows.Item(1).RowHeight = 37

# --- Selection on the new sheet ---
$ws.Range("A1:I1").Select()

# --- Re-activate T2A so it remains the active/visible tab ---
$t2aFinal = $wb.Worksheets.Item("T2A")
$t2aFinal.Activate()

Write-Host "Edit applied successfully"
